# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# Freeze the header row, and wrap the data range A1:U58 in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) ---------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = [string]$cell.Value()
    if ($text.EndsWith("_old")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2404"
    } elseif ($text.EndsWith("_new")) {
        $cell.Value = $text.Substring(0, $text.Length - 4) + "_FV2410"
    }
}

# --- 2. Freeze the header row (pane split after row 1) ------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table (Table1) ---------------
$dataRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
